$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 2634.0833
$ws.Range("J33").Value = 3122.2
$ws.Range("L33").Value = 3122.2
$ws.Range("N33").Value = -3580.2
# Row 43
$ws.Range("H43").Value = 2202031.2
$ws.Range("J43").Value = 2202031.2
$ws.Range("L43").Value = 2202031.2
$ws.Range("N43").Value = -2202169.2
# Row 138
$ws.Range("H138").Value = 125002340
$ws.Range("I138").Value = 2629.3333
$ws.Range("J138").Value = 200002160
$ws.Range("K138").Value = 7887.999899999999
$ws.Range("L138").Value = 600006480
$ws.Range("M138").Value = -2747.999899999999
$ws.Range("N138").Value = -600016760
# Row 141
$ws.Range("H141").Value = 1218.8276
$ws.Range("I141").Value = 1079.1111
$ws.Range("K141").Value = 3237.3333
$ws.Range("M141").Value = 1942.6667

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 11768207
$ws.Range("J32").Value = 4795.3
$ws.Range("L32").Value = 4795.3
$ws.Range("N32").Value = -5369.3
# Row 61
$ws.Range("H61").Value = 23813076
$ws.Range("I61").Value = 32260376
$ws.Range("J61").Value = 7048.5454
$ws.Range("K61").Value = 32260376
$ws.Range("L61").Value = 7048.5454
$ws.Range("M61").Value = -32260164
$ws.Range("N61").Value = -7472.5454
# Row 74
$ws.Range("H74").Value = 21765750
$ws.Range("I74").Value = 25672022
$ws.Range("K74").Value = 25672022
$ws.Range("M74").Value = -25671148
# Row 77
$ws.Range("H77").Value = 21765750
$ws.Range("I77").Value = 25672022
$ws.Range("K77").Value = 128360110
$ws.Range("M77").Value = -128355742
# Row 88
$ws.Range("H88").Value = 52503
$ws.Range("I88").Value = 100006
$ws.Range("J88").Value = 5000
$ws.Range("K88").Value = 100006
$ws.Range("L88").Value = 5000
$ws.Range("M88").Value = -99600
$ws.Range("N88").Value = -5812
# Row 91
$ws.Range("H91").Value = 52503
$ws.Range("I91").Value = 100006
$ws.Range("J91").Value = 5000
$ws.Range("K91").Value = 100006
$ws.Range("L91").Value = 5000
$ws.Range("M91").Value = -98602
$ws.Range("N91").Value = -7808
# Row 132
$ws.Range("H132").Value = 22277610
$ws.Range("I132").Value = 1630.1818
$ws.Range("J132").Value = 83536550
$ws.Range("K132").Value = 4890.5454
$ws.Range("L132").Value = 250609650
$ws.Range("M132").Value = -2360.5454
$ws.Range("N132").Value = -250614710
# Row 136
$ws.Range("H136").Value = 23813076
$ws.Range("I136").Value = 32260376
$ws.Range("J136").Value = 7048.5454
$ws.Range("K136").Value = 96781128
$ws.Range("L136").Value = 21145.6362
$ws.Range("M136").Value = -96778578
$ws.Range("N136").Value = -26245.6362

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 105023100
$ws.Range("I4").Value = 6278750
$ws.Range("J4").Value = 500000500
$ws.Range("K4").Value = 6278750
$ws.Range("L4").Value = 500000500
$ws.Range("M4").Value = -6278638
$ws.Range("N4").Value = -500000724
# Row 6
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()
# Row 8
$ws.Range("H8").Value = 790
$ws.Range("J8").Value = 935
$ws.Range("L8").Value = 935
$ws.Range("N8").Value = -1215
# Row 31
$ws.Range("H31").Value = 25005500
$ws.Range("I31").Value = 3038.4285
$ws.Range("J31").Value = 52639796
$ws.Range("K31").Value = 3038.4285
$ws.Range("L31").Value = 52639796
$ws.Range("M31").Value = -2743.4285
$ws.Range("N31").Value = -52640386
# Row 34
$ws.Range("H34").Value = 25005500
$ws.Range("I34").Value = 3038.4285
$ws.Range("J34").Value = 52639796
$ws.Range("K34").Value = 3038.4285
$ws.Range("L34").Value = 52639796
$ws.Range("M34").Value = -2836.4285
$ws.Range("N34").Value = -52640200
# Row 105
$ws.Range("H105").Value = 7986.926
$ws.Range("J105").Value = 16479.572
$ws.Range("L105").Value = 16479.572
$ws.Range("N105").Value = -19973.572
# Row 141
$ws.Range("H141").Value = 105613.6
$ws.Range("J141").Value = 123267
$ws.Range("L141").Value = 123267
$ws.Range("N141").Value = -133627

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 54589124
$ws.Range("I4").Value = 85755100
$ws.Range("K4").Value = 257265300
$ws.Range("M4").Value = -257265188
# Row 92
$ws.Range("H92").Value = 852.4
$ws.Range("I92").Value = 997.5
$ws.Range("J92").Value = 755.6667
$ws.Range("K92").Value = 2992.5
$ws.Range("L92").Value = 2267.0001
$ws.Range("M92").Value = -1744.5
$ws.Range("N92").Value = -4763.0001
# Row 107
$ws.Range("H107").Value = 602.125
$ws.Range("J107").Value = 636.3
$ws.Range("L107").Value = 1908.9
$ws.Range("N107").Value = -5748.9
# Row 113
$ws.Range("H113").Value = 2163.8125
$ws.Range("J113").Value = 2726
$ws.Range("L113").Value = 8178
$ws.Range("N113").Value = -12518
# Row 128
$ws.Range("H128").Value = 120000
$ws.Range("I128").Value = 120000
$ws.Range("K128").Value = 360000
$ws.Range("M128").Value = -355020

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5260.35
$ws.Range("I70").Value = 4894.6875
$ws.Range("J70").Value = 6723
$ws.Range("K70").Value = 4894.6875
$ws.Range("L70").Value = 6723
$ws.Range("M70").Value = -4624.6875
$ws.Range("N70").Value = -7263
# Row 73
$ws.Range("H73").Value = 5260.35
$ws.Range("I73").Value = 4894.6875
$ws.Range("J73").Value = 6723
$ws.Range("K73").Value = 4894.6875
$ws.Range("L73").Value = 6723
$ws.Range("M73").Value = -3958.6875
$ws.Range("N73").Value = -8595
# Row 102
$ws.Range("H102").Value = 2801.2222
$ws.Range("I102").Value = 1703.1666
$ws.Range("K102").Value = 1703.1666
$ws.Range("M102").Value = -81.16660000000002
# Row 132
$ws.Range("H132").Value = 1594.35
$ws.Range("I132").Value = 1558.0588
$ws.Range("K132").Value = 4674.1764
$ws.Range("M132").Value = -2144.1764

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5665.3335
$ws.Range("I40").Value = 4874.125
$ws.Range("K40").Value = 4874.125
$ws.Range("M40").Value = -4738.125
# Row 93
$ws.Range("H93").Value = 1967.3636
$ws.Range("I93").Value = 1720.4286
$ws.Range("J93").Value = 2399.5
$ws.Range("K93").Value = 1720.4286
$ws.Range("L93").Value = 2399.5
$ws.Range("M93").Value = -472.4286
$ws.Range("N93").Value = -4895.5
# Row 100
$ws.Range("H100").Value = 4111.9
$ws.Range("I100").Value = 3320.7273
$ws.Range("J100").Value = 5078.8887
$ws.Range("K100").Value = 3320.7273
$ws.Range("L100").Value = 5078.8887
$ws.Range("M100").Value = -2779.7273
$ws.Range("N100").Value = -6160.8887
# Row 122
$ws.Range("H122").Value = 11566.333
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 11566.333
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 34698.999
$ws.Range("N122").Value = -39598.999
$ws.Range("M122").ClearContents()
# Row 136
$ws.Range("H136").Value = 3124.5173
$ws.Range("I136").Value = 3037.0454
$ws.Range("J136").Value = 3399.4285
$ws.Range("K136").Value = 9111.136200000001
$ws.Range("L136").Value = 10198.2855
$ws.Range("M136").Value = -6561.136200000001
$ws.Range("N136").Value = -15298.2855

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2359.8
$ws.Range("I81").Value = 800
$ws.Range("J81").Value = 2749.75
$ws.Range("K81").Value = 1600
$ws.Range("L81").Value = 5499.5
$ws.Range("M81").Value = -539
$ws.Range("N81").Value = -7621.5
# Row 84
$ws.Range("H84").Value = 2359.8
$ws.Range("I84").Value = 800
$ws.Range("J84").Value = 2749.75
$ws.Range("K84").Value = 8000
$ws.Range("L84").Value = 27497.5
$ws.Range("M84").Value = -2696
$ws.Range("N84").Value = -38105.5
# Row 126
$ws.Range("H126").Value = 10218.625
$ws.Range("I126").Value = 10535.571
$ws.Range("K126").Value = 31606.713
$ws.Range("M126").Value = -29136.713
# Row 132
$ws.Range("H132").Value = 2311.6667
$ws.Range("I132").Value = 2035.6666
$ws.Range("K132").Value = 6106.9998
$ws.Range("M132").Value = -3576.9998
